$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old A2 value (AUTO_ORG_VPMVD) and set the new data
$ws.Range("A2").ClearContents()

# A1 stays "OrgName" (unchanged value, just shared-string index moved)
$ws.Range("A1").Value = "OrgName"

# New B1 header "SubOrgName", with same style as A1 (s="1")
$ws.Range("B1").Value = "SubOrgName"

# New B2 value
$ws.Range("B2").Value = "AUTO_SUB_ORG_TJLKG"

# Update selection to A3
$ws.Range("A3").Select()
